$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "328.41", "-0.08%") are stored as text, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '328.41'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.08%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.80%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.577'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.59%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08083'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.35%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.911'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.21%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9511'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.68%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.555'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1184'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.34%'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.08%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09755'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.65%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04475'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '6.34%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1067'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.27%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001277'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.23%'
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04204'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-4.26%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005846'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.385'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-4.97%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.305'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.58%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3480'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.07%'
$ws.Range("B20").Value = 'MCDex'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.16'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '15.80%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1416'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '3.84%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001245'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.04%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004347'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.74%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001190'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.91%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-1.24%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02699'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '1.27%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05555'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-0.20%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007553'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.54%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1408'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.44%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007977'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-18.58%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002015'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-5.29%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008395'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-12.76%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007150'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.33%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.09%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.004380'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '25.68%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002268'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.08%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-1.09%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-1.09%'
